$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'244.91"
$ws.Range("D3").Value = "'25.03"
$ws.Range("D6").Value = "'6.520"
$ws.Range("D7").Value = "'2.955"
$ws.Range("D8").Value = "'0.8118"
$ws.Range("D9").Value = "'0.8344"
$ws.Range("D10").Value = "'0.1334"
$ws.Range("D11").Value = "'0.06954"
$ws.Range("D12").Value = "'0.02837"
$ws.Range("D13").Value = "'0.09409"
$ws.Range("D14").Value = "'0.001533"
$ws.Range("D15").Value = "'0.0005963"
$ws.Range("D16").Value = "'0.006064"
$ws.Range("D17").Value = "'3.505"
$ws.Range("D19").Value = "'0.3185"
$ws.Range("D20").Value = "'0.03167"
$ws.Range("D21").Value = "'0.1318"
$ws.Range("D22").Value = "'3.736"
$ws.Range("D23").Value = "'0.04675"
$ws.Range("D25").Value = "'0.001236"
$ws.Range("D27").Value = "'0.00009687"
$ws.Range("E27").Value = "'26NitroExNTXBestin24h"
$ws.Range("D28").Value = "'0.0001950"
$ws.Range("D40").Value = "'0.03622"
$ws.Range("D41").Value = "'0.006252"
$ws.Range("E41").Value = "'40KickTokenKICK"
$ws.Range("D42").Value = "'0.1053"
$ws.Range("D43").Value = "'0.002719"
$ws.Range("D44").Value = "'0.007384"
$ws.Range("D45").Value = "'0.00005273"
$ws.Range("D48").Value = "'0.002284"
